{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the copyright/footer line that follows it, and the blank paragraph that\n// separates them from the preceding \"LOQ4095: ...\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the paragraph that ends the requirements list.\nconst anchorText = \"LOQ4095: Qu\u00edmica Geral Experimental (Requisito)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// The three paragraphs to delete are the ones immediately following the\n// anchor: a blank paragraph, the \"Ver no Jupiter...\" line, and the\n// \"\u00a9 2020 ...\" copyright line. Delete from the bottom up so earlier\n// indices stay valid.\nconst toDelete = [anchorIndex + 1, anchorIndex + 2, anchorIndex + 3];\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  items[toDelete[i]].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the copyright/footer line that follows it, and the blank paragraph that\n# separates them from the preceding \"LOQ4095: ...\" requirement line.\n$d = $word.ActiveDocument\n\n$anchorText = \"LOQ4095: Qu\u00edmica Geral Experimental (Requisito)\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($txt -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# The three paragraphs to delete are the ones immediately following the\n# anchor: a blank paragraph, the \"Ver no Jupiter...\" line, and the\n# \"\u00a9 2020 ...\" copyright line. Delete from the bottom up so earlier\n# indices stay valid.\n$d.Paragraphs.Item($anchorIndex + 3).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 2).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n"}
